$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Merge the split run in paragraph 28 ("Potentially, the best
# solution is for t" + "he man to make four trips...") into a single run,
# and remove the _GoBack bookmark that currently sits between them.
# ---------------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Potentially, the best solution is for t*") {
        $targetPara = $p
        break
    }
}

$paraStart = $targetPara.Range.Start
$paraEnd = $targetPara.Range.End - 1  # exclude the paragraph mark

$splitPoint = $paraStart + "Potentially, the best solution is for t".Length

$secondRunRange = $d.Range($splitPoint, $paraEnd)
$secondRunText = $secondRunRange.Text
$secondRunRange.Delete()

$firstRunRange = $d.Range($paraStart, $splitPoint)
$firstRunRange.InsertAfter($secondRunText)

# ---------------------------------------------------------------------------
# Step 2: Insert the new "4. Evaluate Each Solution" and "5. Choose a
# Solution and Develop a Plan to Implement It" paragraphs right before the
# final (empty) paragraph of the document.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$newContentXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t>4. Evaluate Each Solution:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t>The proposed solution meets the goal and sub-goal. The man transports each item across the river safely and he did not leave any item behind that could potentially cause harm to another item.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t>5. Choose a Solution and Develop a Plan to Implement It:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t xml:space="preserve">The solution that meets the needs of keeping each item safe from the other is for the man to take the parrot with him the first trip. By doing this, the cat is left alone with the seed and it is highly unlikely that the cat will eat the seed. Then, he must go back and get the seed the second trip and </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>when dropping the seed off on the other side of the bank, take the parrot back to the original side again. This ensures that the parrot is not left alone to eat the seed. Once there, drop the parrot off and bring the cat with him to the goal side. This once again leaves the parrot alone on one side and leaves the cat with the seed and no potential harm to any item. Finally, the man makes one more trip to get the parrot and brings it back to the goal bank side and the man and all three items are now safely at their destination.</w:t>
  </w:r>
</w:p>
'@

$insertPoint.InsertXML($newContentXml)

# ---------------------------------------------------------------------------
# Step 3: Re-create the _GoBack bookmark, now collapsed inside the very
# last (empty) paragraph of the document, matching the post-edit location
# Word leaves it in after the most recent typing.
# ---------------------------------------------------------------------------
$docEnd = $d.Content.End
$bookmarkRange = $d.Range($docEnd - 2, $docEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
